# FingerLakes Farms_Triphammer_2025-07-25.xlsx
#
# Update Quantity (column C) and Total Cost (column E) for rows 4-8.
# Cost Per (column D) is unchanged; Total Cost is recomputed as
# Quantity * Cost Per, with Quantity revised down to 1 for these rows.
#
# The values are stored as text in the sheet (not numbers), so entering
# a numeric-looking string via Value would normally get auto-converted
# to a real number by Excel. Using a leading apostrophe forces Excel to
# keep it as text, and resetting the cell Style back to "Normal"
# afterwards clears the transient quote-prefix formatting flag that
# the apostrophe entry leaves behind, so no cell formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $r = $ws.Range($range)
    $r.Value = "'" + $value
    $r.Style = "Normal"
}

Set-TextValue "C4" "1"
Set-TextValue "E4" "9.25"

Set-TextValue "C5" "1"
Set-TextValue "E5" "23.75"

Set-TextValue "C6" "1"
Set-TextValue "E6" "13.00"

Set-TextValue "C7" "1"
Set-TextValue "E7" "13.00"

Set-TextValue "C8" "1"
Set-TextValue "E8" "9.25"
